$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.281.45'
$ws.Range('E2').Value = '  +2.55%  '
$ws.Range('D3').Value = '1.609.99'
$ws.Range('E3').Value = '  +1.22%  '
$ws.Range('E4').Value = '  -0.57%  '
$ws.Range('D5').Value = "'213.06"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.10%  '
$ws.Range('E6').Value = '  -0.57%  '
$ws.Range('E7').Value = '  +0.63%  '
$ws.Range('E8').Value = '  +1.55%  '
$ws.Range('E9').Value = '  +1.70%  '
$ws.Range('D10').Value = "'18.19"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.17%  '
$ws.Range('D11').Value = "'0.0817"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.34%  '
$ws.Range('D12').Value = '1.833.93'
$ws.Range('E12').Value = '  +1.26%  '
$ws.Range('D13').Value = '1.609.71'
$ws.Range('E13').Value = '  +0.95%  '
$ws.Range('E14').Value = '  -0.50%  '
$ws.Range('E15').Value = '  +0.76%  '
$ws.Range('D16').Value = '26.259.50'
$ws.Range('E16').Value = '  +2.50%  '
$ws.Range('D17').Value = "'60.83"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.77%  '
$ws.Range('D18').Value = '0.0₃0729'
$ws.Range('E18').Value = '  +2.12%  '
$ws.Range('E19').Value = '  -0.63%  '
$ws.Range('D20').Value = "'199.18"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.59%  '
$ws.Range('D21').Value = "'4.25"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.92%  '
$ws.Range('D22').Value = "'9.41"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.84%  '
$ws.Range('E23').Value = '  +1.63%  '
$ws.Range('D24').Value = "'0.131"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.36%  '
$ws.Range('D25').Value = "'142.97"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.66%  '
$ws.Range('D26').Value = "'1.75"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.24%  '
$ws.Range('E27').Value = '  -0.60%  '
$ws.Range('D28').Value = "'15.19"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.52%  '
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('E30').Value = '  -0.50%  '
$ws.Range('E31').Value = '  +1.77%  '
$ws.Range('E32').Value = '  +2.24%  '
$ws.Range('D33').Value = "'3.02"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.58%  '
$ws.Range('E34').Value = '  +1.78%  '
$ws.Range('E35').Value = '  -1.61%  '
$ws.Range('D36').Value = '1.109.95'
$ws.Range('E36').Value = '  +1.76%  '
$ws.Range('B37').Value = 'PaxDollar'
$ws.Range('C37').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D37').Value = "'1.00"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.68%  '
$ws.Range('E38').Value = '  +1.02%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = "'2.33"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.84%  '
$ws.Range('E40').Value = '  +1.53%  '
$ws.Range('D41').Value = "'0.788"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.31%  '
$ws.Range('D42').Value = "'0.782"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.84%  '
$ws.Range('D43').Value = '1.746.54'
$ws.Range('E43').Value = '  +1.26%  '
$ws.Range('E44').Value = '  +1.55%  '
$ws.Range('D45').Value = "'92.71"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.75%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').Value = "'1.55"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +9.34%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = "'53.81"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.52%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = "'0.0509"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = "'0.410"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('B50').Value = 'USDD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D50').Value = "'1.00"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.38%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = "'7.32"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.99%  '
